$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" (Overview!G2) is the same shared string as
# "Correspond Handoff Datetime" for de-de (de-de!H2); update both so the
# regenerated report timestamp stays in sync everywhere it is shown.
$wsOverview.Range("G2").Value = "2016-08-25 23:03:17"
$wsDeDe.Range("H2").Value = "2016-08-25 23:03:17"

# zh-cn handoff/handback timestamps
$wsZhCn.Range("H2").Value = "2016-08-25 23:03:12"
$wsZhCn.Range("K2").Value = "2016-08-25 23:03:30"

# de-de handback timestamp
$wsDeDe.Range("K2").Value = "2016-08-25 23:03:37"
